# Apply the "Added list of Russian segment modules" edit to the
# ISS Mimic Magnet Mapping workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix up two module names that gained more specific labels ---------
# Row 14 (Node 2 / Harmony, Starboard) attaches to Columbus -> "COF / Columbus"
$ws.Range("F14").Value = "COF / Columbus"
# Row 25's module name "Columbus" -> "COF /Columbus" (module-name variant, no space)
$ws.Range("A25").Value = "COF /Columbus"
# Row 18 (Node 3 / Tranquility, Forward) attaches to Leonardo -> "PMM / Leonardo"
$ws.Range("F18").Value = "PMM / Leonardo"
# Row 32's module name "Leonardo" -> "PMM / Leonardo"
$ws.Range("A32").Value = "PMM / Leonardo"

# --- Append the new Russian segment module list (rows 45-56) ----------
$ws.Range("A45").Value = "FGB / Zarya"
$ws.Range("B45").Value = "Forward"
$ws.Range("F45").Value = "PMA-1"

$ws.Range("A46").Value = "FGB / Zarya"
$ws.Range("B46").Value = "Nadir"
$ws.Range("F46").Value = "MRM-1 / Rassvet"

$ws.Range("A47").Value = "FGB / Zarya"
$ws.Range("B47").Value = "Aft"
$ws.Range("F47").Value = "Zvezda SM"

$ws.Range("A48").Value = "Zvezda SM"
$ws.Range("B48").Value = "Forward"
$ws.Range("F48").Value = "FGB / Zarya"

$ws.Range("A49").Value = "Zvezda SM"
$ws.Range("B49").Value = "Zenith"
$ws.Range("F49").Value = "MRM-2 / Poisk"

$ws.Range("A50").Value = "Zvezda SM"
$ws.Range("B50").Value = "Nadir"
$ws.Range("F50").Value = "MLM / Nauka"

$ws.Range("A51").Value = "Zvezda SM"
$ws.Range("B51").Value = "Aft"

$ws.Range("A52").Value = "MRM-2 / Poisk"
$ws.Range("B52").Value = "Nadir"
$ws.Range("F52").Value = "Zvezda SM"

$ws.Range("A53").Value = "MRM-2 / Poisk"
$ws.Range("B53").Value = "Zenith"

$ws.Range("A54").Value = "MLM / Nauka"
$ws.Range("B54").Value = "Zenith"
$ws.Range("F54").Value = "Zvezda SM"

$ws.Range("A55").Value = "MLM / Nauka"
$ws.Range("B55").Value = "Nadir"

$ws.Range("A56").Value = "MLM / Nauka"
$ws.Range("B56").Value = "Forward (?)"

# --- Sheet view bookkeeping to match the saved state -------------------
$ws.Range("F56").Select()
